$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "added report message to Mass_Balance_Cor": the E18:E24 mass-balance
# values were reported in the wrong units (x1000 too large) - rescale them.
$ws.Range("E18").Value = 12015 / 1000
$ws.Range("E19").Value = 16984 / 1000
$ws.Range("E20").Value = 2054 / 1000
$ws.Range("E21").Value = 19811 / 1000
$ws.Range("E22").Value = 17541 / 1000
$ws.Range("E23").Value = 14524 / 1000
$ws.Range("E24").Value = 11212 / 1000

# Localize the built-in "Normal" cell style name.
$wb.Styles.Item("Normal").Name = "Normální"

# "separated Start for testing": move the viewport / active selection
# down to the test-start area of the sheet.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E26").Select() | Out-Null
